# "ui buttons and first build test"
#
# - a brand new row 30 is appended with a new time-log entry
# - row 29 (E29) loses its "last row" fill marker (style 8 -> style 7, the
#   same border-only look already used by rows 26-28) and the new row 30
#   picks up that same border/date/duration formatting
# - the U/V summary formulas recalc automatically once the new D30 value
#   exists
# - the view is scrolled/zoomed to a different cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 30 content first, so the SUM(D3:D100)/SUM(M3:M100) style ranges
# pick the new value up during recalculation.
$ws.Cells.Item(30, 1).Value = 45361
$ws.Cells.Item(30, 2).Formula = "=16+42/60"
$ws.Cells.Item(30, 3).Formula = "=17+10/60"
$ws.Cells.Item(30, 4).Formula = "=C30-B30"
$ws.Cells.Item(30, 5).Value = "UI buttons working and first build test"

# Copy row 29's formatting down into the new row 30 (borders, date/duration
# number formats, etc.)
$ws.Range("A29:E29").Copy() | Out-Null
$ws.Range("A30:E30").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# Row 29 drops the distinct "last row" fill flag now that row 30 exists.
$ws.Range("E29").Interior.Pattern = -4142

# View state: zoomed out a bit and scrolled to T10.
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("T10").Select() | Out-Null
